$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Columns("M").Delete()
$ws.Range("M1").Select()
